$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting from the row above (D8) so the new date cell
# picks up the same style (short date format) instead of Excel inventing
# a brand new number format.
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Fill in row 9 with the new week's data (Miercoles 05/07/2025)
$ws.Range("D9").Value = 45843
$ws.Range("E9").Value = 419
$ws.Range("F9").Value = 367
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 300
$ws.Range("J9").Value = "N/A"

# Update the sheet view: scroll back to the left edge (A1) and move the
# active selection further down, below the newly added row.
$ws.Range("A1").Select()
$ws.Range("J19").Select()
